$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: insert a block of "flow content" (runs / bookmarks / etc, no
# paragraph marks) at the very start of a paragraph without causing the
# engine to merge it with the paragraph's existing runs and without
# accidentally splitting the paragraph in two.  Inserting an OOXML
# fragment wrapped in a single <w:p> exactly at a paragraph's Start
# position merges the fragment's children into that paragraph instead of
# creating a new one.
# ---------------------------------------------------------------------------
function Insert-FlowXmlAtStart {
    param($doc, [int]$pos, [string]$innerXml)

    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p>' + $innerXml + '</w:p></w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'

    $rng = $doc.Range($pos, $pos)
    $rng.InsertXML($pkg)
}

# ---------------------------------------------------------------------------
# Step 0: remove the old "_GoBack" bookmark that used to sit at the very end
# of the document.  It gets re-created further below, right after the new
# "story" run.  Doing the removal first (before any text is touched) avoids
# ever having two "_GoBack" bookmarks in the document at once.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# Step 1: rebuild the paragraph that talks about "seven stories" / "a story"
# ---------------------------------------------------------------------------

# Locate the paragraph by its (unique) text.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*you will read seven stories*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find target paragraph"
}

$pStart = $target.Range.Start
$pEnd = $target.Range.End

$rPr = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:sz w:val="24"/></w:rPr>'

$newContent =
    '<w:r>' + $rPr + '<w:t xml:space="preserve">This study is on what people think another person knows. Specifically, you will read </w:t></w:r>' +
    '<w:r>' + $rPr + '<w:t>a</w:t></w:r>' +
    '<w:r>' + $rPr + '<w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r>' + $rPr + '<w:t>story</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '<w:r>' + $rPr + '<w:t xml:space="preserve"> about </w:t></w:r>' +
    '<w:r>' + $rPr + '<w:t>people and what they think to be true</w:t></w:r>' +
    '<w:r>' + $rPr + '<w:t xml:space="preserve">, then you will </w:t></w:r>' +
    '<w:r>' + $rPr + '<w:t>answer a few questions about these stories</w:t></w:r>' +
    '<w:r>' + $rPr + '<w:t xml:space="preserve">, followed by a </w:t></w:r>' +
    '<w:r>' + $rPr + '<w:t xml:space="preserve">very </w:t></w:r>' +
    '<w:r>' + $rPr + '<w:t xml:space="preserve">brief questionnaire about yourself. </w:t></w:r>'

# Insert the brand-new, fully-formed run/bookmark sequence right before the
# paragraph's existing (old) content.
Insert-FlowXmlAtStart $d $pStart $newContent

# The old content got pushed after what we just inserted; figure out where
# it now lives and delete it (the paragraph mark itself must stay put).
$newLen = $target.Range.Text.Length
$oldTextLen = ($pEnd - $pStart) - 1   # exclude the paragraph mark
$oldStart = $pStart + $newLen - $oldTextLen - 1
$oldRange = $d.Range($oldStart, $pStart + $newLen - 1)
$oldRange.Delete()

